# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text in A1 -------------
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 10.92 = 44377.73 pesos
✅ 44377.73 pesos = 10.87 = 961.56 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$ws1.Range("A1").Value = $newText

# --- tasas: refresh the auto-updated rate cells -----------------------------
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 91.59999999999999
$ws2.Range("O10").Value = 4065
$ws2.Range("N12").Value = 4083
$ws2.Range("O12").Value = 88.46899999999999
